$d = $word.ActiveDocument

$pairs = @(
    @{old="205×5=1025"; new="361×5=1805"},
    @{old="196×5=980"; new="996×8=7968"},
    @{old="239×9=2151"; new="937×6=5622"},
    @{old="320×8=2560"; new="808×9=7272"},
    @{old="265×6=1590"; new="865×6=5190"},
    @{old="912×8=7296"; new="353×6=2118"},
    @{old="868×8=6944"; new="950×4=3800"},
    @{old="202×5=1010"; new="644×6=3864"},
    @{old="377×3=1131"; new="159×2=318"},
    @{old="430×6=2580"; new="266×4=1064"},
    @{old="402×3=1206"; new="972×5=4860"},
    @{old="908×9=8172"; new="464×5=2320"},
    @{old="512×6=3072"; new="761×7=5327"},
    @{old="310×9=2790"; new="410×2=820"},
    @{old="764×5=3820"; new="448×2=896"},
    @{old="906×7=6342"; new="519×8=4152"},
    @{old="201×7=1407"; new="857×5=4285"},
    @{old="360×9=3240"; new="109×5=545"},
    @{old="502×4=2008"; new="968×3=2904"},
    @{old="396×2=792"; new="163×6=978"},
    @{old="428×4=1712"; new="740×8=5920"},
    @{old="862×3=2586"; new="409×6=2454"},
    @{old="661×9=5949"; new="277×5=1385"},
    @{old="565×6=3390"; new="202×4=808"},
    @{old="372×8=2976"; new="664×9=5976"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
